$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 2.8
$ws.Range("H4").Value = 3.05
$ws.Range("I4").Value = 2.47
$ws.Range("V4").Value = 2
$ws.Range("W4").Value = 9
$ws.Range("AB4").Value = 29
$ws.Range("AC4").Value = 9.25
$ws.Range("AD4").Value = 5.9
$ws.Range("AE4").Value = 12.5
$ws.Range("AH4").Value = 8.25
$ws.Range("AI4").Value = 12.5
$ws.Range("AK4").Value = 28
$ws.Range("AM4").Value = 28
$ws.Range("AP4").Value = 19
$ws.Range("AQ4").Value = 65
$ws.Range("AU4").Value = 6.3
$ws.Range("AV4").Value = 45
$ws.Range("AW4").Value = 4.5
$ws.Range("AX4").Value = 13
$ws.Range("AY4").Value = 18.5
$ws.Range("AZ4").Value = 55
$ws.Range("BA4").Value = 75
